$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5 (existing rows 5-15 shift down to 6-16)
$ws.Rows.Item(5).EntireRow.Insert()

# Populate the newly inserted row 5 with the new "Carahue" record
$ws.Cells.Item(5, 1).Value = 11
$ws.Cells.Item(5, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(5, 3).Value = "Bíobío"
$ws.Cells.Item(5, 4).Value = 44519
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = 100112022
$ws.Cells.Item(5, 7).Value = "Arveja Verde"
$ws.Cells.Item(5, 8).Value = "Perfection"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 240
$ws.Cells.Item(5, 11).Value = 17000
$ws.Cells.Item(5, 12).Value = 18000
$ws.Cells.Item(5, 13).Value = 17583
$ws.Cells.Item(5, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(5, 15).Value = "Carahue"
$ws.Cells.Item(5, 16).Value = 703
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"
